# Resolve conflicts then merge branch 'master' into feature/laglobal.
#
# The "Effective end date" column (D) is dropped from the overtime import
# template. Deleting the entire column shifts "Effective End Time" (old
# column E) left into D, removes the now-unused "Effective end date" shared
# string, and shrinks the sheet's used range from A1:E1 to A1:D1 - all in
# one coherent operation, matching how this was produced in Excel/Calc.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep the sheet's default column width in sync (cosmetic, best-effort).
$ws.StandardWidth = 9.01171875

# Remove column D ("Effective end date"); column E ("Effective End Time")
# and everything to its right shifts one column to the left.
$ws.Columns.Item(4).Delete()

# The active selection moves from the old C7 to D7 after the column removal.
$ws.Range("D7").Select()
